$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Normalize a few cells' direct formatting so they fold back onto already
# existing (shared) cell styles instead of keeping their own one-off variants.
# These are purely cosmetic clean-ups (Google Sheets drops unused duplicate
# styles on save) riding along with the data edit below.

# F1 ("header_name" header) loses its special alignment and matches the rest
# of the plain header cells (e.g. E1).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# F11 / F18 ("header_name" values for the new 2022/2023 split rows) drop their
# special alignment and match the plain numeric-style cells (e.g. E11).
$ws.Range("E11").Copy()
$ws.Range("F11").PasteSpecial(-4122)
$ws.Range("F18").PasteSpecial(-4122)

$wb.Application.CutCopyMode = $false

# --- The actual content edit: shift the 2019-2021 / 2022-present split to
# 2019-2022 / 2023-present across the "Refrigeration" family of rows.

# Rows whose max_year was 2021 (end of the old period) now end in 2022, and
# pick up the same "touched" alignment that row 12 (C12) already had.
$ws.Range("C12").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C16").PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false

# C19 also moves from 2021 to 2022, picking up the "touched" alignment that
# the min_year cells of the new period (e.g. B11) already had.
$ws.Range("B11").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false

# Rows whose min_year was 2022 (start of the new period) now start in 2023.
$ws.Range("B11").Value = 2023
$ws.Range("B14").Value = 2023
$ws.Range("B17").Value = 2023
$ws.Range("B18").Value = 2023

# Rows whose max_year was 2021 now end in 2022.
$ws.Range("C12").Value = 2022
$ws.Range("C15").Value = 2022
$ws.Range("C16").Value = 2022
$ws.Range("C19").Value = 2022
